# Apply "Todays Changes To Flows and Control Relationships"
#
# Adds two new annotation cells in column H of the "Claim Filing" sheet,
# widens the new column, and moves the active selection to G14 to match
# the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Claim Filing")

# New call-out text next to the NMR / BPTW decision rows
$ws.Range("H4").Value = "What if NMR is NOT PC ?"
$ws.Range("H8").Value = "What if BPTW is NOT PC ?"

# Give the new column enough room for the text (matches authored width)
$ws.Columns.Item(8).ColumnWidth = 26.5

# Leave the selection where the author left it after editing
$ws.Activate()
$ws.Range("G14").Select()
